$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the style/formatting used by the other header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), rows 2-33
$data = @{
    2  = @(8, 8)
    3  = @(2, 3)
    4  = @(1, 1)
    5  = @(6, 7)
    6  = @(3, 4)
    7  = @(9, 9)
    8  = @(7, 7)
    9  = @(7, 7)
    10 = @(6, 6)
    11 = @(8, 8)
    12 = @(7, 7)
    13 = @(8, 8)
    14 = @(7, 7)
    15 = @(9, 9)
    16 = @(8, 9)
    17 = @(7, 7)
    18 = @(6, 7)
    19 = @(8, 8)
    20 = @(9, 9)
    21 = @(5, 5)
    22 = @(9, 9)
    23 = @(9, 9)
    24 = @(8, 8)
    25 = @(5, 5)
    26 = @(4, 4)
    27 = @(6, 6)
    28 = @(6, 6)
    29 = @(6, 6)
    30 = @(5, 6)
    31 = @(6, 6)
    32 = @(7, 7)
    33 = @(2, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
